{"js": "// Refresh the \"Ice and Fire\" review copy: new title, reworked pros/cons\n// bullets, and an updated closing summary, per the commit \"Added many\n// more features\".\n\nconst replacements = [\n  // Title heading + the later bold recap line (same text, two locations).\n  [\n    \"Play Ice and Fire Free Slot - A Unique Dragon-Themed Game\",\n    \"Play Ice and Fire for Free\",\n  ],\n  // \"What we like\" bullets\n  [\n    \"Unique dual-screen gameplay mechanics\",\n    \"Unique gameplay with dual 5x5 grids\",\n  ],\n  [\n    \"Special symbols and bonus features for more winning opportunities\",\n    \"Special symbols trigger bonus features\",\n  ],\n  [\n    \"Visually stunning graphics and immersive sound effects\",\n    \"Visually appealing dual-screen design\",\n  ],\n  [\n    \"Dragon-themed game with both Eastern and Western cultures\",\n    \"Engaging dragon-themed experience\",\n  ],\n  // \"What we don't like\" bullets\n  [\n    \"The game may be too complex for inexperienced players\",\n    \"Limited interaction between the grids\",\n  ],\n  [\n    \"The lack of interaction between the two playing fields may make the game feel repetitive\",\n    \"Requires consecutive combinations to access special mode\",\n  ],\n  // Closing italic summary line\n  [\n    \"Read our review of Ice and Fire, a unique dragon-themed online slot game. Play for free and experience dual-screen gameplay and special symbols.\",\n    \"Read a review of Ice and Fire, a unique dragon-themed online slot game. Play for free and experience the dual-screen design.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Ice and Fire\" content refresh described in the commit.\n# Uses Word's Find/Replace (wdReplaceAll) across the whole document body\n# for each old -> new text pair, matching the unified diff exactly.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n# Title (appears twice: the Heading1 title and the later bold recap line)\nReplace-All \"Play Ice and Fire Free Slot - A Unique Dragon-Themed Game\" \"Play Ice and Fire for Free\"\n\n# \"What we like\" bullets\nReplace-All \"Unique dual-screen gameplay mechanics\" \"Unique gameplay with dual 5x5 grids\"\nReplace-All \"Special symbols and bonus features for more winning opportunities\" \"Special symbols trigger bonus features\"\nReplace-All \"Visually stunning graphics and immersive sound effects\" \"Visually appealing dual-screen design\"\nReplace-All \"Dragon-themed game with both Eastern and Western cultures\" \"Engaging dragon-themed experience\"\n\n# \"What we don't like\" bullets\nReplace-All \"The game may be too complex for inexperienced players\" \"Limited interaction between the grids\"\nReplace-All \"The lack of interaction between the two playing fields may make the game feel repetitive\" \"Requires consecutive combinations to access special mode\"\n\n# Closing meta summary (bold title line already handled above) + italic description\nReplace-All \"Read our review of Ice and Fire, a unique dragon-themed online slot game. Play for free and experience dual-screen gameplay and special symbols.\" \"Read a review of Ice and Fire, a unique dragon-themed online slot game. Play for free and experience the dual-screen design.\"\n"}
